$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Avg. Marks" -> "Total Marks"
$ws.Range("C1").Value = "Total Marks"

# Row 2 (Labib): Total Marks = 89*12
$ws.Range("C2").Formula = "=89*12"

# Row 3 (Namir): Total Marks = 88*12
$ws.Range("C3").Formula = "=88*12"

# Row 4 (Shafin): Total Marks = 88*12
$ws.Range("C4").Formula = "=88*12"

# Row 5 (Tamim): clear the Total Marks / Comment cells entirely
$ws.Range("C5:D5").ClearContents()

# Row 6 (Zamim): Total Marks = 84*12
$ws.Range("C6").Formula = "=84*12"

# Row 7 (Shafiq): clear Roll, Total Marks = 79*12
$ws.Range("B7").ClearContents()
$ws.Range("C7").Formula = "=79*12"

# Row 8 (Asir): Total Marks = 78*12
$ws.Range("C8").Formula = "=78*12"

# Row 9 (Nasir): Total Marks = 76*12
$ws.Range("C9").Formula = "=76*12"

# Row 10 (Basit): Total Marks = 74*12
$ws.Range("C10").Formula = "=74*12"

# Row 11 (Mihir): Total Marks = 72*12
$ws.Range("C11").Formula = "=72*12"

# Column widths: C widened (closest attainable to 15.5703125 given character-width
# quantization), D newly widened to 16
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666

# Selection moves to F18
[void]$ws.Range("F18").Select()
